# Added generation of impulse and ramp wind entries in the DLC_List sheet.

$wb = $excel.ActiveWorkbook

$wsConfig = $wb.Worksheets.Item("config")
$wsDlc = $wb.Worksheets.Item("DLC_List")

# New DLC rows: "impulse" and "ramp" wind definitions
$wsDlc.Range("A6").Value = "impulse"
$wsDlc.Range("B6").Value = "IMP:+1/5"
$wsDlc.Range("C6").Value = "[v_r-2 v_r v_r+2]"
$wsDlc.Range("F6").Value = "0"
$wsDlc.Range("G6").Value = "0"
$wsDlc.Range("H6").Value = "150"
$wsDlc.Range("I6").Value = "50"

$wsDlc.Range("A7").Value = "ramp"
$wsDlc.Range("B7").Value = "RMP:50"
$wsDlc.Range("C7").Value = "3"
$wsDlc.Range("F7").Value = "0"
$wsDlc.Range("G7").Value = "0"
$wsDlc.Range("H7").Value = "1100"
$wsDlc.Range("I7").Value = "50"

# Update the remembered cursor/selection position on each sheet
$wsConfig.Activate()
$wsConfig.Range("B14").Select()

$wsDlc.Activate()
$wsDlc.Range("H8").Select()
